$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39:123 down to 40:124.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record's data.
$ws.Cells.Item(39, 1).Value = 4
$ws.Cells.Item(39, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value = "Los Lagos"
$ws.Cells.Item(39, 4).Value = 44536
$ws.Cells.Item(39, 5).Value = 10
$ws.Cells.Item(39, 6).Value = "Fruta"
$ws.Cells.Item(39, 7).Value = 100108
$ws.Cells.Item(39, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(39, 9).Value = 100108002
$ws.Cells.Item(39, 10).Value = "Mango"
$ws.Cells.Item(39, 11).Value = "Sin especificar"
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value = 150
$ws.Cells.Item(39, 14).Value = 5000
$ws.Cells.Item(39, 15).Value = 8000
$ws.Cells.Item(39, 16).Value = 6833
$ws.Cells.Item(39, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(39, 18).Value = "Perú"
$ws.Cells.Item(39, 19).Value = 1708
$ws.Cells.Item(39, 20).Value = 4
